# Fruta / hortaliza, semanal
# Insert two new weekly-report rows (new survey date 2022-01-24 / serial 44585)
# at the top of the "Platano" price block, pushing the existing rows
# (formerly 465-491) down to 467-493.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before row 465 (shifts old 465..491 down to 467..493)
$ws.Range("A465:A466").EntireRow.Insert()

# New row 465: "Pintón" quality
$ws.Range("A465").Value = 7
$ws.Range("B465").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C465").Value = "Ñuble"
$ws.Range("D465").Value = 44585
$ws.Range("E465").Value = 16
$ws.Range("F465").Value = "Fruta"
$ws.Range("G465").Value = 100108
$ws.Range("H465").Value = "Tropicales y subtropicales"
$ws.Range("I465").Value = 100108006
$ws.Range("J465").Value = "Plátano"
$ws.Range("K465").Value = "Sin especificar"
$ws.Range("L465").Value = "Pintón"
$ws.Range("M465").Value = 100
$ws.Range("N465").Value = 13000
$ws.Range("O465").Value = 13000
$ws.Range("P465").Value = 13000
$ws.Range("Q465").Value = "`$/caja 20 kilos"
$ws.Range("R465").Value = "Ecuador"
$ws.Range("S465").Value = 650
$ws.Range("T465").Value = 20

# New row 466: "Primera Pintón" quality
$ws.Range("A466").Value = 7
$ws.Range("B466").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C466").Value = "Ñuble"
$ws.Range("D466").Value = 44585
$ws.Range("E466").Value = 16
$ws.Range("F466").Value = "Fruta"
$ws.Range("G466").Value = 100108
$ws.Range("H466").Value = "Tropicales y subtropicales"
$ws.Range("I466").Value = 100108006
$ws.Range("J466").Value = "Plátano"
$ws.Range("K466").Value = "Sin especificar"
$ws.Range("L466").Value = "Primera Pintón"
$ws.Range("M466").Value = 200
$ws.Range("N466").Value = 14000
$ws.Range("O466").Value = 15000
$ws.Range("P466").Value = 14500
$ws.Range("Q466").Value = "`$/caja 20 kilos"
$ws.Range("R466").Value = "Ecuador"
$ws.Range("S466").Value = 725
$ws.Range("T466").Value = 20
